$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.875.59"
$ws.Range("E2").Value = "  -2.55%  "

$ws.Range("D3").Value = "3.874.43"
$ws.Range("E3").Value = "  -2.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.671"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.753"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("E10").Value = "  +3.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000323"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.97%  "

$ws.Range("D14").Value = "4.489.49"
$ws.Range("E14").Value = "  -2.59%  "

$ws.Range("D15").Value = "3.852.20"
$ws.Range("E15").Value = "  -2.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.33%  "

$ws.Range("E19").Value = "  -2.03%  "

$ws.Range("D20").Value = "70.734.02"
$ws.Range("E20").Value = "  -2.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "436.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.51%  "

$ws.Range("E24").Value = "  -4.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "48.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "70.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.49%  "

$ws.Range("E35").Value = "  -4.23%  "

$ws.Range("D36").Value = "0.0₃0979"
$ws.Range("E36").Value = "  +11.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "628.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.424"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +29.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.144"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  -2.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0469"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.19%  "

$ws.Range("E47").Value = "  -3.36%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -15.56%  "

$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.84%  "

$ws.Range("D50").Value = "2.826.13"
$ws.Range("E50").Value = "  +1.04%  "

$ws.Range("E51").Value = "  +0.06%  "
